$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.801.09"
$ws.Range("E2").Value = "  +1.43%  "

$ws.Range("D3").Value = "2.089.49"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.34%  "

$ws.Range("E9").Value = "  +1.59%  "

$ws.Range("E10").Value = "  +2.35%  "

$ws.Range("E11").Value = "  +2.90%  "

$ws.Range("D12").Value = "2.384.72"
$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.43"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.58%  "

$ws.Range("E15").Value = "  -1.43%  "

$ws.Range("E16").Value = "  +2.39%  "

$ws.Range("D17").Value = "2.090.36"
$ws.Range("E17").Value = "  +1.33%  "

$ws.Range("D18").Value = "37.683.54"
$ws.Range("E18").Value = "  +1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.15%  "

$ws.Range("D21").Value = "0.0₃0821"
$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.15"
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  -1.61%  "

$ws.Range("E25").Value = "  -0.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.142"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.29%  "

$ws.Range("E28").Value = "  +1.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("E32").Value = "  +3.65%  "

$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.80%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.44%  "

$ws.Range("E41").Value = "  -0.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.96%  "

$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("D44").Value = "1.453.24"
$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("E45").Value = "  -0.24%  "

$ws.Range("E46").Value = "  +3.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.58%  "

$ws.Range("E48").Value = "  +4.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.12%  "

$ws.Range("E50").Value = "  +2.11%  "

$ws.Range("D51").Value = "2.280.24"
$ws.Range("E51").Value = "  +1.20%  "
